$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4620.115
$ws.Range("J17").Value = 5153.609
$ws.Range("L17").Value = 15460.827
$ws.Range("N17").Value = -15796.827
$ws.Range("H64").Value = 4160.0435
$ws.Range("I64").Value = 3587.5
$ws.Range("J64").Value = 4600.4614
$ws.Range("K64").Value = 3587.5
$ws.Range("L64").Value = 4600.4614
$ws.Range("M64").Value = -3339.5
$ws.Range("N64").Value = -5096.4614
$ws.Range("H67").Value = 4160.0435
$ws.Range("I67").Value = 3587.5
$ws.Range("J67").Value = 4600.4614
$ws.Range("K67").Value = 3587.5
$ws.Range("L67").Value = 4600.4614
$ws.Range("M67").Value = -2729.5
$ws.Range("N67").Value = -6316.4614
$ws.Range("H76").Value = 3468.6458
$ws.Range("I76").Value = 3285.7144
$ws.Range("K76").Value = 3285.7144
$ws.Range("M76").Value = -2970.7144
$ws.Range("H79").Value = 3468.6458
$ws.Range("I79").Value = 3285.7144
$ws.Range("K79").Value = 3285.7144
$ws.Range("M79").Value = -2193.7144
$ws.Range("H113").Value = 2941.15
$ws.Range("I113").Value = 2508.3333
$ws.Range("J113").Value = 3295.2727
$ws.Range("K113").Value = 2508.3333
$ws.Range("L113").Value = 3295.2727
$ws.Range("M113").Value = 745.6667000000002
$ws.Range("N113").Value = -9803.2727
$ws.Range("H129").Value = 1072.2609
$ws.Range("I129").Value = 312.5
$ws.Range("J129").Value = 1144.619
$ws.Range("K129").Value = 937.5
$ws.Range("L129").Value = 3433.857
$ws.Range("M129").Value = 4062.5
$ws.Range("N129").Value = -13433.857
$ws.Range("H137").Value = 753349.3
$ws.Range("I137").Value = 2451.7693
$ws.Range("J137").Value = 2528198
$ws.Range("K137").Value = 7355.3079
$ws.Range("L137").Value = 7584594
$ws.Range("M137").Value = -4805.3079
$ws.Range("N137").Value = -7589694
$ws.Range("H138").Value = 4884.345
$ws.Range("I138").Value = 1157.4736
$ws.Range("J138").Value = 7963.0654
$ws.Range("K138").Value = 3472.4208
$ws.Range("L138").Value = 23889.1962
$ws.Range("M138").Value = 1667.5792
$ws.Range("N138").Value = -34169.19620000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20917.807
$ws.Range("I32").Value = 22271.25
$ws.Range("K32").Value = 22271.25
$ws.Range("M32").Value = -21984.25
$ws.Range("H43").Value = 15277
$ws.Range("J43").Value = 15277
$ws.Range("L43").Value = 15277
$ws.Range("N43").Value = -15903
$ws.Range("H45").Value = 1579.0646
$ws.Range("I45").Value = 1575.4445
$ws.Range("J45").Value = 1603.5
$ws.Range("K45").Value = 1575.4445
$ws.Range("L45").Value = 1603.5
$ws.Range("M45").Value = -1198.4445
$ws.Range("N45").Value = -2357.5
$ws.Range("H63").Value = 2625.625
$ws.Range("I63").Value = 2500.8333
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 2500.8333
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -1814.8333
$ws.Range("N63").Value = -4372
$ws.Range("H66").Value = 2625.625
$ws.Range("I66").Value = 2500.8333
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 12504.1665
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -9072.166499999999
$ws.Range("N66").Value = -21864
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = $null
$ws.Range("N120").Value = 0
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 217
$ws.Range("I22").Value = 217
$ws.Range("K22").Value = 217
$ws.Range("M22").Value = 133
$ws.Range("H31").Value = 4167.46
$ws.Range("I31").Value = 4660.9
$ws.Range("J31").Value = 3427.3
$ws.Range("K31").Value = 4660.9
$ws.Range("L31").Value = 3427.3
$ws.Range("M31").Value = -4365.9
$ws.Range("N31").Value = -4017.3
$ws.Range("H34").Value = 4167.46
$ws.Range("I34").Value = 4660.9
$ws.Range("J34").Value = 3427.3
$ws.Range("K34").Value = 4660.9
$ws.Range("L34").Value = 3427.3
$ws.Range("M34").Value = -4458.9
$ws.Range("N34").Value = -3831.3
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 15257.571
$ws.Range("J68").Value = 17667.166
$ws.Range("L68").Value = 53001.49800000001
$ws.Range("N68").Value = -54623.49800000001
$ws.Range("H70").Value = 3037.5
$ws.Range("I70").Value = 1433.3334
$ws.Range("K70").Value = 4300.0002
$ws.Range("M70").Value = -3985.0002
$ws.Range("H71").Value = 15257.571
$ws.Range("J71").Value = 17667.166
$ws.Range("L71").Value = 159004.494
$ws.Range("N71").Value = -167116.494
$ws.Range("H73").Value = 3037.5
$ws.Range("I73").Value = 1433.3334
$ws.Range("K73").Value = 4300.0002
$ws.Range("M73").Value = -3208.0002
$ws.Range("H121").Value = 299
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = $null
$ws.Range("H122").Value = 976
$ws.Range("I122").Value = 575
$ws.Range("J122").Value = 997.67566
$ws.Range("K122").Value = 5175
$ws.Range("L122").Value = 8979.08094
$ws.Range("M122").Value = -2725
$ws.Range("N122").Value = -13879.08094
$ws.Range("H131").Value = 32400.645
$ws.Range("J131").Value = 54559.555
$ws.Range("L131").Value = 163678.665
$ws.Range("N131").Value = -173758.665
$ws.Range("H137").Value = 20650.4
$ws.Range("I137").Value = 1699.0834
$ws.Range("J137").Value = 33284.61
$ws.Range("K137").Value = 5097.2502
$ws.Range("L137").Value = 99853.83
$ws.Range("M137").Value = 2.749799999999595
$ws.Range("N137").Value = -110053.83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 35725.832
$ws.Range("J64").Value = 35725.832
$ws.Range("L64").Value = 35725.832
$ws.Range("N64").Value = -36221.832
$ws.Range("H67").Value = 35725.832
$ws.Range("J67").Value = 35725.832
$ws.Range("L67").Value = 35725.832
$ws.Range("N67").Value = -37441.832
$ws.Range("H80").Value = 15000.625
$ws.Range("J80").Value = 8200
$ws.Range("L80").Value = 8200
$ws.Range("N80").Value = -10196
$ws.Range("H83").Value = 15000.625
$ws.Range("J83").Value = 8200
$ws.Range("L83").Value = 41000
$ws.Range("N83").Value = -50984
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = $null
$ws.Range("M17").Value = $null
$ws.Range("N17").Value = 0
$ws.Range("H22").Value = 714.52
$ws.Range("I22").Value = 567.5
$ws.Range("J22").Value = 850.2308
$ws.Range("K22").Value = 567.5
$ws.Range("L22").Value = 850.2308
$ws.Range("M22").Value = -272.5
$ws.Range("N22").Value = -1440.2308
$ws.Range("H27").Value = 714.52
$ws.Range("I27").Value = 567.5
$ws.Range("J27").Value = 850.2308
$ws.Range("K27").Value = 567.5
$ws.Range("L27").Value = 850.2308
$ws.Range("M27").Value = -460.5
$ws.Range("N27").Value = -1064.2308
$ws.Range("H47").Value = 11065
$ws.Range("J47").Value = 11065
$ws.Range("L47").Value = 11065
$ws.Range("N47").Value = -12045
$ws.Range("H52").Value = 11065
$ws.Range("J52").Value = 11065
$ws.Range("L52").Value = 11065
$ws.Range("N52").Value = -11531
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = $null
$ws.Range("N27").Value = 0
$ws.Range("H62").Value = 3890.3572
$ws.Range("J62").Value = 3860.4546
$ws.Range("L62").Value = 3860.4546
$ws.Range("N62").Value = -5108.4546
$ws.Range("H65").Value = 3890.3572
$ws.Range("J65").Value = 3860.4546
$ws.Range("L65").Value = 19302.273
$ws.Range("N65").Value = -25542.273
